# Applies the "Updated Sequence Diagrams and Ui Diagram" edit:
#  - rename deletePerson(p) -> deleteTask(t)
#  - rename AddressBookChangedEvent -> TaskManagerChangedEvent (x2)
#  - rename handleAddresssBookChangedEvent -> handleTaskManagerChangedEvent (x2)
#  - refresh the "datetimeFigureOut" date placeholders (10/16/2016 -> 3/24/17)
#  - (best effort) register the two slide guides shown in the deck

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Sequence-diagram label edits on slide 1.
#    Each textbox keeps its original run split, so we patch the characters in
#    place (via TextRange.Characters) instead of overwriting the whole range,
#    which would otherwise merge/re-split the runs and could disturb the
#    per-run formatting (color, err flag, etc).
# ---------------------------------------------------------------------------

function Set-RunText($shape, [int]$start, [string]$oldText, [string]$newText) {
    $len = $oldText.Length
    $rng = $shape.TextFrame.TextRange.Characters($start, $len)
    if ($rng.Text -eq $oldText) {
        $rng.Text = $newText
    }
}

$slide = $p.Slides.Item(1)

# "deletePerson(p)" -> "deleteTask(t)"
$sh = $slide.Shapes.Item("TextBox 28")
Set-RunText $sh 1 "deletePerson" "deleteTask"
Set-RunText $sh 11 "(p)" "(t)"

# "post(AddressBookChangedEvent)" -> "post(TaskManagerChangedEvent)"
foreach ($name in @("TextBox 32", "TextBox 61")) {
    $sh = $slide.Shapes.Item($name)
    Set-RunText $sh 6 "AddressBookChangedEvent" "TaskManagerChangedEvent"
}

# "handleAddresssBookChangedEvent()" -> "handleTaskManagerChangedEvent()"
foreach ($name in @("TextBox 73", "TextBox 49")) {
    $sh = $slide.Shapes.Item($name)
    Set-RunText $sh 1 "handleAddresssBookChangedEvent" "handleTaskManagerChangedEvent"
}

# ---------------------------------------------------------------------------
# 2. Refresh the "date updated automatically" placeholders that live on the
#    slide master and every slide layout (and, where the host allows writes,
#    the notes master) from 10/16/2016 to 3/24/17.
# ---------------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "10/16/2016") {
                $tr.Text = "3/24/17"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

if ($p.HasNotesMaster) {
    Update-DatePlaceholder $p.NotesMaster.Shapes
}

# ---------------------------------------------------------------------------
# 3. Best-effort: make sure the two slide guides (horizontal @1488,
#    vertical @2880) used throughout the deck are registered on the
#    presentation.
# ---------------------------------------------------------------------------

try {
    $p.Guides.Add(2, 1488) | Out-Null
    $p.Guides.Add(1, 2880) | Out-Null
} catch {
}
